$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# Remove rows 12-18 from Neg_Change (symbols no longer present after the refresh)
$ws1.Range("A12:I18").EntireRow.Delete()

# Update Neg_Change rows 2-11 with the refreshed market data
$ws1.Cells.Item(2,1).Value = "ASIANPAINT"
$ws1.Cells.Item(2,2).Value = 2906.4
$ws1.Cells.Item(2,3).Value = 2926.9
$ws1.Cells.Item(2,4).Value = 2880
$ws1.Cells.Item(2,5).Value = 2885.4
$ws1.Cells.Item(2,6).Value = 1230683
$ws1.Cells.Item(2,7).Value = 2494076
$ws1.Cells.Item(2,8).Value = -0.5065575387438074
$ws1.Cells.Item(2,9).Value = "ASIANPAINT"

$ws1.Cells.Item(3,1).Value = "AUROPHARMA"
$ws1.Cells.Item(3,2).Value = 1222.5
$ws1.Cells.Item(3,3).Value = 1241.1
$ws1.Cells.Item(3,4).Value = 1218.3
$ws1.Cells.Item(3,5).Value = 1240
$ws1.Cells.Item(3,6).Value = 988088
$ws1.Cells.Item(3,7).Value = 2055973
$ws1.Cells.Item(3,8).Value = -0.5194061400611778
$ws1.Cells.Item(3,9).Value = "AUROPHARMA"

$ws1.Cells.Item(4,1).Value = "TIINDIA"
$ws1.Cells.Item(4,2).Value = 3067.6
$ws1.Cells.Item(4,3).Value = 3106.6
$ws1.Cells.Item(4,4).Value = 3053
$ws1.Cells.Item(4,5).Value = 3096
$ws1.Cells.Item(4,6).Value = 107644
$ws1.Cells.Item(4,7).Value = 214533
$ws1.Cells.Item(4,8).Value = -0.4982403639533312
$ws1.Cells.Item(4,9).Value = "TIINDIA"

$ws1.Cells.Item(5,1).Value = "ALKEM"
$ws1.Cells.Item(5,2).Value = 5750
$ws1.Cells.Item(5,3).Value = 5750
$ws1.Cells.Item(5,4).Value = 5631.5
$ws1.Cells.Item(5,5).Value = 5725
$ws1.Cells.Item(5,6).Value = 176721
$ws1.Cells.Item(5,7).Value = 383633
$ws1.Cells.Item(5,8).Value = -0.5393488047170082
$ws1.Cells.Item(5,9).Value = "ALKEM"

$ws1.Cells.Item(6,1).Value = "SJVN"
$ws1.Cells.Item(6,2).Value = 83.69
$ws1.Cells.Item(6,3).Value = 84.4
$ws1.Cells.Item(6,4).Value = 83.25
$ws1.Cells.Item(6,5).Value = 83.35
$ws1.Cells.Item(6,6).Value = 1892839
$ws1.Cells.Item(6,7).Value = 3737752
$ws1.Cells.Item(6,8).Value = -0.4935889272482498
$ws1.Cells.Item(6,9).Value = "SJVN"

$ws1.Cells.Item(7,1).Value = "PIIND"
$ws1.Cells.Item(7,2).Value = 3576.3
$ws1.Cells.Item(7,3).Value = 3582
$ws1.Cells.Item(7,4).Value = 3528
$ws1.Cells.Item(7,5).Value = 3545
$ws1.Cells.Item(7,6).Value = 84301
$ws1.Cells.Item(7,7).Value = 167191
$ws1.Cells.Item(7,8).Value = -0.4957802752540508
$ws1.Cells.Item(7,9).Value = "PIIND"

$ws1.Cells.Item(8,1).Value = "ATGL"
$ws1.Cells.Item(8,2).Value = 630
$ws1.Cells.Item(8,3).Value = 630
$ws1.Cells.Item(8,4).Value = 620.6
$ws1.Cells.Item(8,5).Value = 621.8
$ws1.Cells.Item(8,6).Value = 387009
$ws1.Cells.Item(8,7).Value = 931132
$ws1.Cells.Item(8,8).Value = -0.5843672003539777
$ws1.Cells.Item(8,9).Value = "ATGL"

$ws1.Cells.Item(9,1).Value = "MANKIND"
$ws1.Cells.Item(9,2).Value = 2280
$ws1.Cells.Item(9,3).Value = 2280
$ws1.Cells.Item(9,4).Value = 2236
$ws1.Cells.Item(9,5).Value = 2244
$ws1.Cells.Item(9,6).Value = 292698
$ws1.Cells.Item(9,7).Value = 606234
$ws1.Cells.Item(9,8).Value = -0.5171864329615297
$ws1.Cells.Item(9,9).Value = "MANKIND"

$ws1.Cells.Item(10,1).Value = "NBCC"
$ws1.Cells.Item(10,2).Value = 114.7
$ws1.Cells.Item(10,3).Value = 117.5
$ws1.Cells.Item(10,4).Value = 113.7
$ws1.Cells.Item(10,5).Value = 116.95
$ws1.Cells.Item(10,6).Value = 15262875
$ws1.Cells.Item(10,7).Value = 37966431
$ws1.Cells.Item(10,8).Value = -0.5979902614496474
$ws1.Cells.Item(10,9).Value = "NBCC"

$ws1.Cells.Item(11,1).Value = "RBLBANK"
$ws1.Cells.Item(11,2).Value = 321
$ws1.Cells.Item(11,3).Value = 321
$ws1.Cells.Item(11,4).Value = 315.95
$ws1.Cells.Item(11,5).Value = 317.15
$ws1.Cells.Item(11,6).Value = 2615820
$ws1.Cells.Item(11,7).Value = 5156226
$ws1.Cells.Item(11,8).Value = -0.4926870932344704
$ws1.Cells.Item(11,9).Value = "RBLBANK"

# Update Pos_Change rows 2-7 with the refreshed market data
$ws2.Cells.Item(2,1).Value = "LODHA"
$ws2.Cells.Item(2,2).Value = 1212.4
$ws2.Cells.Item(2,3).Value = 1223
$ws2.Cells.Item(2,4).Value = 1210.5
$ws2.Cells.Item(2,5).Value = 1221.3
$ws2.Cells.Item(2,6).Value = 781591
$ws2.Cells.Item(2,7).Value = 541800
$ws2.Cells.Item(2,8).Value = 0.4425821336286452
$ws2.Cells.Item(2,9).Value = "LODHA"

$ws2.Cells.Item(3,1).Value = "ZYDUSLIFE"
$ws2.Cells.Item(3,2).Value = 948
$ws2.Cells.Item(3,3).Value = 949.95
$ws2.Cells.Item(3,4).Value = 932.15
$ws2.Cells.Item(3,5).Value = 935
$ws2.Cells.Item(3,6).Value = 873274
$ws2.Cells.Item(3,7).Value = 592731
$ws2.Cells.Item(3,8).Value = 0.4733057660220235
$ws2.Cells.Item(3,9).Value = "ZYDUSLIFE"

$ws2.Cells.Item(4,1).Value = "MARICO"
$ws2.Cells.Item(4,2).Value = 753.5
$ws2.Cells.Item(4,3).Value = 764.65
$ws2.Cells.Item(4,4).Value = 750.2
$ws2.Cells.Item(4,5).Value = 758.8
$ws2.Cells.Item(4,6).Value = 3726013
$ws2.Cells.Item(4,7).Value = 2572321
$ws2.Cells.Item(4,8).Value = 0.448502344769568
$ws2.Cells.Item(4,9).Value = "MARICO"

$ws2.Cells.Item(5,1).Value = "HINDPETRO"
$ws2.Cells.Item(5,2).Value = 483
$ws2.Cells.Item(5,3).Value = 492.2
$ws2.Cells.Item(5,4).Value = 482.4
$ws2.Cells.Item(5,5).Value = 486
$ws2.Cells.Item(5,6).Value = 3478922
$ws2.Cells.Item(5,7).Value = 2456049
$ws2.Cells.Item(5,8).Value = 0.4164709254579204
$ws2.Cells.Item(5,9).Value = "HINDPETRO"

$ws2.Cells.Item(6,1).Value = "COLPAL"
$ws2.Cells.Item(6,2).Value = 2172.9
$ws2.Cells.Item(6,3).Value = 2194.6
$ws2.Cells.Item(6,4).Value = 2172.9
$ws2.Cells.Item(6,5).Value = 2185.7
$ws2.Cells.Item(6,6).Value = 135290
$ws2.Cells.Item(6,7).Value = 96584
$ws2.Cells.Item(6,8).Value = 0.4007496065600928
$ws2.Cells.Item(6,9).Value = "COLPAL"

$ws2.Cells.Item(7,1).Value = "BANDHANBNK"
$ws2.Cells.Item(7,2).Value = 154.99
$ws2.Cells.Item(7,3).Value = 157.75
$ws2.Cells.Item(7,4).Value = 154.55
$ws2.Cells.Item(7,5).Value = 155
$ws2.Cells.Item(7,6).Value = 5810479
$ws2.Cells.Item(7,7).Value = 4141559
$ws2.Cells.Item(7,8).Value = 0.4029690268809402
$ws2.Cells.Item(7,9).Value = "BANDHANBNK"

